$d = $word.ActiveDocument

# The first paragraph holds the "**ID__...__ID**" merge-field placeholder
# (currently split across two runs: the marker text, then a trailing
# space-only run).
$p1 = $d.Paragraphs.Item(1)

# --- Replace the paragraph's text (everything but the trailing paragraph
#     mark) with the updated placeholder. This collapses both existing runs
#     into a single run and drops the trailing space. ---
$pRange = $p1.Range
$textOnly = $d.Range($pRange.Start, $pRange.End - 1)
$textOnly.Text = "**ID__AFFARS_AFMC_PGI_5307_104__ID**"

# --- Paragraph formatting: add a paragraph border (5-twip padding on all
#     sides) and widen the left indent from 120 to 225 twips. ---
$p1 = $d.Paragraphs.Item(1)
$p1.Range.ParagraphFormat.LeftIndent = 11.25
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
